# New weekly record: insert a row at row 50 (pushing the existing rows
# 50-92 down to 51-93) and populate it with the new day's price data for
# Ciruela (Black Amber / Primera) at Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(50).Insert()

$ws.Range("A50").Value = 9
$ws.Range("B50").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C50").Value = "Metropolitana"
$ws.Range("D50").Value = 44603
$ws.Range("E50").Value = 13
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100103
$ws.Range("H50").Value = "Frutos de hueso (carozo)"
$ws.Range("I50").Value = 100103002
$ws.Range("J50").Value = "Ciruela"
$ws.Range("K50").Value = "Black Amber"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 350
$ws.Range("N50").Value = 9000
$ws.Range("O50").Value = 9000
$ws.Range("P50").Value = 9000
$ws.Range("Q50").Value = "$/caja 15 kilos granel"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 600
$ws.Range("T50").Value = 15
